$wb = $excel.ActiveWorkbook

# ----- "Logs" sheet: append two new mail-log rows (19, 20) -----
$logs = $wb.Worksheets.Item("Logs")

$logs.Cells.Item(19, 1).Value = "Klacht over levering"
$logs.Cells.Item(19, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(19, 3).Value = "Ik ben niet tevreden over mijn bestelling. Ik hoor graag hoe jullie dit oplossen."
$logs.Cells.Item(19, 4).Value = "Klacht"
$logs.Cells.Item(19, 6).Value = "2025-06-17 16:00:10"
$logs.Cells.Item(19, 7).Value = "Nee"

$logs.Cells.Item(20, 1).Value = "Vragen over samenwerking"
$logs.Cells.Item(20, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item(20, 3).Value = "Kunnen we samenwerken aan een nieuw project?"
$logs.Cells.Item(20, 4).Value = "Overig"
$logs.Cells.Item(20, 6).Value = "2025-06-17 16:30:23"
$logs.Cells.Item(20, 7).Value = "Nee"

# Extend the conditional-formatting ranges so the new rows are covered,
# keeping them in sync with the grown data range (D2:D18 -> D2:D20,
# G2:G18 -> G2:G20).
$logs.Range("D2:D18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D20"))
$logs.Range("G2:G18").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G20"))

# ----- "Dashboard" sheet: refresh the category summary counts -----
$dash = $wb.Worksheets.Item("Dashboard")

# "Overig" gained one more mail (the new "Vragen over samenwerking" row).
$dash.Cells.Item(3, 2).Value = 6

# "Klacht" gained one more mail (the new "Klacht over levering" row) and
# moves up to row 5; "Bestelling" shifts down to row 6, count unchanged.
$dash.Cells.Item(5, 1).Value = "Klacht"
$dash.Cells.Item(5, 2).Value = 2
$dash.Cells.Item(6, 1).Value = "Bestelling"
$dash.Cells.Item(6, 2).Value = 1
